$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194, pushing the existing rows 194:240 down to 195:241
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new record
$ws.Cells.Item(194, 1).Value = 3
$ws.Cells.Item(194, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(194, 3).Value = 'Coquimbo'
$ws.Cells.Item(194, 4).Value = 44508
$ws.Cells.Item(194, 5).Value = 5
$ws.Cells.Item(194, 6).Value = 100112031
$ws.Cells.Item(194, 7).Value = 'Poroto verde'
$ws.Cells.Item(194, 8).Value = 'Magnum'
$ws.Cells.Item(194, 9).Value = 'Primera'
$ws.Cells.Item(194, 10).Value = 73
$ws.Cells.Item(194, 11).Value = 40000
$ws.Cells.Item(194, 12).Value = 41000
$ws.Cells.Item(194, 13).Value = 40521
$ws.Cells.Item(194, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(194, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(194, 16).Value = 1621
$ws.Cells.Item(194, 17).Value = 25
$ws.Cells.Item(194, 18).Value = 'Hortaliza'
